$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold font, thin border, center/top alignment) from
# the existing "IP" header cell (H1) onto the two new header cells so they
# match the look of the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-45
$iVals = @(7,7,10,6,4,6,6,9,7,7,6,8,9,7,1,6,6,6,6,10,9,6,7,9,9,7,7,8,6,7,7,8,9,6,8,5,1,5,5,6,4,6,1,3)
$jVals = @(7,7,10,6,6,6,7,9,7,7,6,8,9,7,3,7,6,6,6,10,9,6,7,9,9,7,7,8,6,7,7,8,9,6,8,6,1,5,5,6,4,6,2,3)

for ($idx = 0; $idx -lt 44; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
